$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 52499.5
$ws.Range("J18").Value = 52499.5
$ws.Range("L18").Value = 52499.5
$ws.Range("N18").Value = -53067.5
$ws.Range("H46").Value = 2389.375
$ws.Range("J46").Value = 1699.6666
$ws.Range("L46").Value = 5098.9998
$ws.Range("N46").Value = -5336.9998
$ws.Range("H60").Value = 2389.375
$ws.Range("J60").Value = 1699.6666
$ws.Range("L60").Value = 5098.9998
$ws.Range("N60").Value = -6066.9998
$ws.Range("H82").Value = 18036.5
$ws.Range("I82").Value = 699.6667
$ws.Range("K82").Value = 2099.0001
$ws.Range("M82").Value = -1693.0001
$ws.Range("H85").Value = 18036.5
$ws.Range("I85").Value = 699.6667
$ws.Range("K85").Value = 2099.0001
$ws.Range("M85").Value = -695.0001000000002
$ws.Range("H132").Value = 2595.2942
$ws.Range("I132").Value = 2062.6538
$ws.Range("K132").Value = 6187.9614
$ws.Range("M132").Value = -3657.9614
$ws.Range("H133").Value = 67998.2
$ws.Range("J133").Value = 67998.2
$ws.Range("L133").Value = 67998.2
$ws.Range("N133").Value = -78118.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 767
$ws.Range("I37").Value = 767
$ws.Range("K37").Value = 767
$ws.Range("M37").Value = -494
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0
$ws.Range("H88").Value = 2984.6667
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 3469.3333
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 3469.3333
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -4281.3333
$ws.Range("H91").Value = 2984.6667
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 3469.3333
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 3469.3333
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -6277.3333
$ws.Range("H97").Value = 790.6
$ws.Range("I97").Value = 798
$ws.Range("K97").Value = 798
$ws.Range("M97").Value = -302
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 43060
$ws.Range("J114").Value = 43060
$ws.Range("L114").Value = 43060
$ws.Range("N114").Value = -51738
$ws.Range("H120").Value = 55300
$ws.Range("J120").Value = 55300
$ws.Range("L120").Value = 55300
$ws.Range("N120").Value = -64976
$ws.Range("H123").Value = 113000
$ws.Range("J123").Value = 113000
$ws.Range("L123").Value = 113000
$ws.Range("N123").Value = -122800
$ws.Range("H134").Value = 30895.27
$ws.Range("I134").Value = 1504.7916
$ws.Range("J134").Value = 85154.62
$ws.Range("K134").Value = 4514.3748
$ws.Range("L134").Value = 255463.86
$ws.Range("M134").Value = -1979.3748
$ws.Range("N134").Value = -260533.86
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 50460
$ws.Range("J18").Value = 50460
$ws.Range("L18").Value = 50460
$ws.Range("N18").Value = -50920
$ws.Range("H22").Value = 1160.5454
$ws.Range("I22").Value = 1251.6
$ws.Range("K22").Value = 1251.6
$ws.Range("M22").Value = -901.5999999999999
$ws.Range("H52").Value = 85996
$ws.Range("J52").Value = 85996
$ws.Range("L52").Value = 85996
$ws.Range("N52").Value = -86584
$ws.Range("H98").Value = 55936.332
$ws.Range("J98").Value = 55936.332
$ws.Range("L98").Value = 55936.332
$ws.Range("N98").Value = -60428.332
$ws.Range("H112").Value = 54166.332
$ws.Range("J112").Value = 54166.332
$ws.Range("L112").Value = 54166.332
$ws.Range("N112").Value = -57120.332
$ws.Range("H115").Value = 40963
$ws.Range("J115").Value = 40963
$ws.Range("L115").Value = 40963
$ws.Range("N115").Value = -43313
$ws.Range("H116").Value = 98595.75
$ws.Range("J116").Value = 98595.75
$ws.Range("L116").Value = 98595.75
$ws.Range("N116").Value = -107773.75
$ws.Range("H117").Value = 85400
$ws.Range("J117").Value = 85400
$ws.Range("L117").Value = 85400
$ws.Range("N117").Value = -94578
$ws.Range("H119").Value = 80059.5
$ws.Range("J119").Value = 80059.5
$ws.Range("L119").Value = 80059.5
$ws.Range("N119").Value = -89735.5
$ws.Range("H127").Value = 79959.75
$ws.Range("J127").Value = 79959.75
$ws.Range("L127").Value = 79959.75
$ws.Range("N127").Value = -89879.75
$ws.Range("H134").Value = 720069.0600000001
$ws.Range("I134").Value = 1251375.2
$ws.Range("K134").Value = 3754125.6
$ws.Range("M134").Value = -3751590.6
$ws.Range("H137").Value = 81995.336
$ws.Range("J137").Value = 81995.336
$ws.Range("L137").Value = 81995.336
$ws.Range("N137").Value = -92195.336
$ws.Range("H138").Value = 77465.664
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 77465.664
$ws.Range("K138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("M138").Value = 77465.664
$ws.Range("N138").Value = -87745.664
$ws.Range("H139").Value = 97379.5
$ws.Range("J139").Value = 97379.5
$ws.Range("L139").Value = 97379.5
$ws.Range("N139").Value = -107659.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20112.7
$ws.Range("I2").Value = 126.90909
$ws.Range("J2").Value = 44539.777
$ws.Range("K2").Value = 761.4545400000001
$ws.Range("L2").Value = 267238.662
$ws.Range("M2").Value = -648.4545400000001
$ws.Range("N2").Value = -267464.662
$ws.Range("H25").Value = 5000100
$ws.Range("I25").Value = 200
$ws.Range("J25").Value = 10000000
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 30000000
$ws.Range("M25").Value = -431
$ws.Range("N25").Value = -30000338
$ws.Range("H30").Value = 5000100
$ws.Range("I30").Value = 200
$ws.Range("J30").Value = 10000000
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 30000000
$ws.Range("M30").Value = -498
$ws.Range("N30").Value = -30000204
$ws.Range("H137").Value = 5971.5
$ws.Range("J137").Value = 3780.8333
$ws.Range("L137").Value = 11342.4999
$ws.Range("N137").Value = -21542.4999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3263
$ws.Range("I80").Value = 3253.2144
$ws.Range("K80").Value = 3253.2144
$ws.Range("M80").Value = -2255.2144
$ws.Range("H83").Value = 3263
$ws.Range("I83").Value = 3253.2144
$ws.Range("K83").Value = 16266.072
$ws.Range("M83").Value = -11274.072
$ws.Range("H122").Value = 1999.5652
$ws.Range("I122").Value = 2017.6316
$ws.Range("J122").Value = 1913.75
$ws.Range("K122").Value = 6052.8948
$ws.Range("L122").Value = 5741.25
$ws.Range("M122").Value = -3602.8948
$ws.Range("N122").Value = -10641.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 92598.5
$ws.Range("J36").Value = 92598.5
$ws.Range("L36").Value = 92598.5
$ws.Range("N36").Value = -93722.5
$ws.Range("H46").Value = 4015.762
$ws.Range("I46").Value = 3784.3572
$ws.Range("J46").Value = 4478.5713
$ws.Range("K46").Value = 3784.3572
$ws.Range("L46").Value = 4478.5713
$ws.Range("M46").Value = -3596.3572
$ws.Range("N46").Value = -4854.5713
$ws.Range("H119").Value = 90990
$ws.Range("J119").Value = 90990
$ws.Range("L119").Value = 90990
$ws.Range("N119").Value = -100666
$ws.Range("H121").Value = 88716.664
$ws.Range("J121").Value = 88716.664
$ws.Range("L121").Value = 88716.664
$ws.Range("N121").Value = -92210.664
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("H125").Value = 99995
$ws.Range("J125").Value = 99995
$ws.Range("L125").Value = 99995
$ws.Range("N125").Value = -109835
$ws.Range("H131").Value = 66666.664
$ws.Range("J131").Value = 66666.664
$ws.Range("L131").Value = 66666.664
$ws.Range("N131").Value = -76746.664
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 253499.5
$ws.Range("J17").Value = 5499.5
$ws.Range("L17").Value = 5499.5
$ws.Range("N17").Value = -5843.5
$ws.Range("H128").Value = 67830
$ws.Range("J128").Value = 67830
$ws.Range("L128").Value = 67830
$ws.Range("N128").Value = -77790
$ws.Range("H129").Value = 92215
$ws.Range("J129").Value = 92215
$ws.Range("L129").Value = 92215
$ws.Range("N129").Value = -102215
$ws.Range("H130").Value = 88489
$ws.Range("J130").Value = 88489
$ws.Range("L130").Value = 88489
$ws.Range("N130").Value = -98529
$ws.Range("H131").Value = 78868
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 78868
$ws.Range("K131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("M131").Value = 78868
$ws.Range("N131").Value = -88948
$ws.Range("H136").Value = 12021.869
$ws.Range("J136").Value = 31673
$ws.Range("L136").Value = 95019
$ws.Range("N136").Value = -100119
